$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text in cell E8 from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active selection on the sheet (E8) as captured in the saved view state
$ws.Activate()
$ws.Range("E8").Select()
